$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: RandomForestRegressor - name unchanged, values updated
$ws.Range("B3").Value = 0.9210965097580663
$ws.Range("C3").Value = 0.9194637266586801
$ws.Range("D3").Value = 0.9000468613481928

# Row 4: GradientBoostingRegressor -> DecisionTreeRegressor
$ws.Range("A4").Value = "DecisionTreeRegressor"
$ws.Range("B4").Value = 0.8229420134181807
$ws.Range("C4").Value = 0.8181393877578844
$ws.Range("D4").Value = 0.5818122877977095

# Row 5: AdaBoostRegressor -> MLPRegressor
$ws.Range("A5").Value = "MLPRegressor"
$ws.Range("B5").Value = 0.8488063417379207
$ws.Range("C5").Value = 0.8431544986279862
$ws.Range("D5").Value = 0.8247663922657517
